# Fixed Swap Column Problems in the Main Chart. Updated to 1.6r.
#
# The last header column (X1) had a typo in its label; correct it.
# "活动那个内容网盘链接" -> "活动内容网盘链接"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Replace("活动那个内容网盘链接", "活动内容网盘链接")

# Move / restore the active selection to the last header cell (X1),
# matching where the author left the cursor after the edit.
$ws.Range("X1").Select()
